$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '64.017.02'
Set-TextValue 'E2' '  +0.83%  '
Set-TextValue 'D3' '3.136.55'
Set-TextValue 'E3' '  +0.61%  '
Set-TextValue 'E4' '  -0.04%  '
Set-TextValue 'D5' '602.62'
Set-TextValue 'E5' '  -0.56%  '
Set-TextValue 'E6' '  -1.21%  '
Set-TextValue 'E7' '  +0.11%  '
Set-TextValue 'D8' '3.131.42'
Set-TextValue 'E8' '  +0.55%  '
Set-TextValue 'D9' '0.522'
Set-TextValue 'E9' '  +0.47%  '
Set-TextValue 'E10' '  -0.48%  '
Set-TextValue 'D11' '5.37'
Set-TextValue 'E11' '  +1.49%  '
Set-TextValue 'E12' '  -0.25%  '
Set-TextValue 'E13' '  +0.64%  '
Set-TextValue 'D14' '35.34'
Set-TextValue 'E14' '  +0.32%  '
Set-TextValue 'D15' '3.658.37'
Set-TextValue 'E15' '  +0.51%  '
Set-TextValue 'E16' '  +2.74%  '
Set-TextValue 'D17' '64.158.46'
Set-TextValue 'E17' '  +0.92%  '
Set-TextValue 'D18' '3.147.11'
Set-TextValue 'E18' '  +0.70%  '
Set-TextValue 'D19' '6.83'
Set-TextValue 'E19' '  +0.35%  '
Set-TextValue 'D20' '479.94'
Set-TextValue 'E20' '  +1.45%  '
Set-TextValue 'D21' '14.65'
Set-TextValue 'E21' '  +1.15%  '
Set-TextValue 'D22' '0.709'
Set-TextValue 'E22' '  +0.04%  '
Set-TextValue 'D23' '7.64'
Set-TextValue 'E23' '  -2.64%  '
Set-TextValue 'D24' '88.09'
Set-TextValue 'E24' '  +5.93%  '
Set-TextValue 'D25' '13.39'
Set-TextValue 'E25' '  -0.91%  '
Set-TextValue 'E26' '  +0.00%  '
Set-TextValue 'D27' '2.74'
Set-TextValue 'E27' '  -1.37%  '
Set-TextValue 'D28' '8.31'
Set-TextValue 'E28' '  -1.42%  '
Set-TextValue 'D29' '7.11'
Set-TextValue 'E29' '  +1.54%  '
Set-TextValue 'D30' '2.07'
Set-TextValue 'E30' '  +0.69%  '
Set-TextValue 'E31' '  -6.43%  '
Set-TextValue 'B32' 'EthereumClassic'
Set-TextValue 'C32' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D32' '27.08'
Set-TextValue 'E32' '  +3.72%  '
Set-TextValue 'B33' 'FirstDigitalUSD'
Set-TextValue 'C33' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D33' '1.00'
Set-TextValue 'E33' '  +0.05%  '
Set-TextValue 'D34' '2.66'
Set-TextValue 'E34' '  -0.21%  '
Set-TextValue 'E35' '  -1.32%  '
Set-TextValue 'D36' '6.04'
Set-TextValue 'E36' '  +2.01%  '
Set-TextValue 'D37' '0.0₃0754'
Set-TextValue 'E37' '  -3.74%  '
Set-TextValue 'D38' '52.62'
Set-TextValue 'E38' '  -0.07%  '
Set-TextValue 'D39' '2.99'
Set-TextValue 'E39' '  +0.50%  '
Set-TextValue 'D40' '439.50'
Set-TextValue 'E40' '  -3.40%  '
Set-TextValue 'D41' '0.0395'
Set-TextValue 'E41' '  +0.95%  '
Set-TextValue 'E42' '  +0.28%  '
Set-TextValue 'D43' '8.27'
Set-TextValue 'E43' '  -0.11%  '
Set-TextValue 'D44' '2.874.53'
Set-TextValue 'E44' '  +1.02%  '
Set-TextValue 'D45' '0.261'
Set-TextValue 'E45' '  -1.43%  '
Set-TextValue 'E46' '  +2.11%  '
Set-TextValue 'E47' '  -2.57%  '
Set-TextValue 'D49' '25.93'
Set-TextValue 'E49' '  -0.63%  '
Set-TextValue 'E50' '  +0.55%  '
Set-TextValue 'D51' '121.87'
Set-TextValue 'E51' '  +2.58%  '
